$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# C2: empty inline string -> numeric value
$ws.Range("C2").Value = 3.536391946991341

# C3: 5.106645561260892 -> 6.025160828942479
$ws.Range("C3").Value = 6.025160828942479

# E3: empty inline string -> numeric value
$ws.Range("E3").Value = 9.522007334472145

# C4: 6.089698944253206 -> 9.409620166348361
$ws.Range("C4").Value = 9.409620166348361
